$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.095.92"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "1.731.95"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'310.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.38%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "'0.4857"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.04%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "'42.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "'0.07291"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'20.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "'5.900"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "1.729.80"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "'6.902"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.72%  "
$ws.Range("D17").Value = "'87.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.73%  "
$ws.Range("D18").Value = "'0.00001040"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'0.9992"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "'16.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").Value = "'5.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "27.138.42"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("D24").Value = "'10.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").Value = "'2.081"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.51%  "
$ws.Range("D26").Value = "'153.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.30%  "
$ws.Range("D27").Value = "'19.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "1.934.63"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").Value = "'2.100"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("D30").Value = "'121.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("D31").Value = "'1.048"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.44%  "
$ws.Range("D32").Value = "'0.09323"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "'3.618"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").Value = "'5.422"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").Value = "'0.02199"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "'0.05930"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("E37").Value = "  -5.87%  "
$ws.Range("D38").Value = "'1.430"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.56%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'4.793"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2005"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("D41").Value = "'0.6016"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.78%  "
$ws.Range("D42").Value = "'0.9982"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'1.097"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.01%  "
$ws.Range("D44").Value = "'7.514"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").Value = "'3.588"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("D47").Value = "'0.5678"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").Value = "'118.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").Value = "'1.850"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.93%  "
$ws.Range("D50").Value = "'1.110"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "'0.06651"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.86%  "
